$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = $ws.Range("I1").Value2
$ws.Range("I1").Value = $ws.Range("J1").Value2
$ws.Range("J1").ClearContents()

$ws.Range("H2").Value = $ws.Range("I2").Value2
$ws.Range("I2").Value = $ws.Range("J2").Value2
$ws.Range("J2").ClearContents()
